{"js": "// Apply the AMOVA covariance table edits:\n//  1. Widen the first table-grid column from 4777 -> 5045 (dxa)\n//  2. \"Among populations within urban/rural groups\" -> \"Among sampling sites within urban/rural groups\"\n//  3. Bump the \"Within ...\" row height from 612 -> 614 (dxa)\n//  4. \"Within populations\" -> \"Within sampling sites\"\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// --- 1. Resize the first column (gridCol 4777 -> 5045 dxa = 238.85pt -> 252.25pt) ---\nconst firstColCell = table.getCell(0, 0);\nfirstColCell.columnWidth = 5045 / 20; // Word.js columnWidth is expressed in points\nawait context.sync();\n\n// --- 2 & 4. Text replacements inside the table cells ---\nconst amongResults = context.document.body.search(\n  \"Among populations within urban/rural groups\",\n  { matchCase: true }\n);\namongResults.load(\"items\");\nawait context.sync();\namongResults.items[0].insertText(\n  \"Among sampling sites within urban/rural groups\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst withinResults = context.document.body.search(\"Within populations\", {\n  matchCase: true,\n});\nwithinResults.load(\"items\");\nawait context.sync();\nwithinResults.items[0].insertText(\n  \"Within sampling sites\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 3. Row height change for the \"Within sampling sites\" row (612 -> 614 dxa) ---\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst withinRow = rows.items[3]; // 0:header, 1:Among urban/rural, 2:Among sampling sites..., 3:Within sampling sites, 4:Total\nwithinRow.preferredHeight = 614 / 20; // preferredHeight is in points\nawait context.sync();\n", "ps1": "# Apply the AMOVA covariance table edits:\n#  1. Widen the first table-grid column from 4777 -> 5045 (dxa)\n#  2. \"Among populations within urban/rural groups\" -> \"Among sampling sites within urban/rural groups\"\n#  3. Bump the \"Within ...\" row height from 612 -> 614 (dxa)\n#  4. \"Within populations\" -> \"Within sampling sites\"\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# --- 1. Resize the first column (gridCol 4777 -> 5045 dxa = 252.25pt) ---\n$table.Columns.Item(1).Width = 5045 / 20\n\n# --- 2. \"Among populations within urban/rural groups\" -> \"Among sampling sites within urban/rural groups\" ---\n$find1 = $d.Content.Find\n$find1.Execute(\"Among populations within urban/rural groups\", $false, $false, $false, $false, $false, $true, 1, $false, \"Among sampling sites within urban/rural groups\", 2)\n\n# --- 4. \"Within populations\" -> \"Within sampling sites\" ---\n$find2 = $d.Content.Find\n$find2.Execute(\"Within populations\", $false, $false, $false, $false, $false, $true, 1, $false, \"Within sampling sites\", 2)\n\n# --- 3. Row height change for the \"Within sampling sites\" row (612 -> 614 dxa) ---\n$table.Rows.Item(4).Height = 614 / 20\n"}
